$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H111").Value = 861.86365
$ws.Range("I111").Value = 485.2353
$ws.Range("J111").Value = 2142.4
$ws.Range("K111").Value = 1455.7059
$ws.Range("L111").Value = 6427.200000000001
$ws.Range("M111").Value = 1611.2941
$ws.Range("N111").Value = -12561.2

$ws.Range("H132").Value = 246711.81
$ws.Range("I132").Value = 253918.83
$ws.Range("J132").Value = 184251
$ws.Range("K132").Value = 761756.49
$ws.Range("L132").Value = 552753
$ws.Range("M132").Value = -759226.49
$ws.Range("N132").Value = -557813

$ws.Range("H135").Value = 1140.1277
$ws.Range("I135").Value = 1011.0455
$ws.Range("J135").Value = 3033.3333
$ws.Range("K135").Value = 9099.4095
$ws.Range("L135").Value = 27299.9997
$ws.Range("M135").Value = -6564.4095
$ws.Range("N135").Value = -32369.9997

$ws.Range("H137").Value = 17544764
$ws.Range("I137").Value = 21739746
$ws.Range("J137").Value = 2110.7273
$ws.Range("K137").Value = 65219238
$ws.Range("L137").Value = 6332.1819
$ws.Range("M137").Value = -65216688
$ws.Range("N137").Value = -11432.1819

$ws.Range("H138").Value = 1050.13
$ws.Range("I138").Value = 462.70587
$ws.Range("J138").Value = 1661.5306
$ws.Range("K138").Value = 1388.11761
$ws.Range("L138").Value = 4984.5918
$ws.Range("M138").Value = 3751.88239
$ws.Range("N138").Value = -15264.5918

$ws.Range("H141").Value = 2168.1775
$ws.Range("I141").Value = 1271.3269
$ws.Range("K141").Value = 3813.9807
$ws.Range("M141").Value = 1366.0193

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2026.8103
$ws.Range("I61").Value = 1621.3721
$ws.Range("J61").Value = 3189.0667
$ws.Range("K61").Value = 1621.3721
$ws.Range("L61").Value = 3189.0667
$ws.Range("M61").Value = -1409.3721
$ws.Range("N61").Value = -3613.0667

$ws.Range("H97").Value = 5285.3335
$ws.Range("I97").Value = 6404.7646
$ws.Range("J97").Value = 527.75
$ws.Range("K97").Value = 6404.7646
$ws.Range("L97").Value = 527.75
$ws.Range("M97").Value = -5908.7646
$ws.Range("N97").Value = -1519.75

$ws.Range("H122").Value = 1448.742
$ws.Range("I122").Value = 1324.875
$ws.Range("J122").Value = 1873.4286
$ws.Range("K122").Value = 3974.625
$ws.Range("L122").Value = 5620.2858
$ws.Range("M122").Value = -1524.625
$ws.Range("N122").Value = -10520.2858

$ws.Range("H136").Value = 2026.8103
$ws.Range("I136").Value = 1621.3721
$ws.Range("J136").Value = 3189.0667
$ws.Range("K136").Value = 4864.1163
$ws.Range("L136").Value = 9567.2001
$ws.Range("M136").Value = -2314.1163
$ws.Range("N136").Value = -14667.2001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1239.8
$ws.Range("I20").Value = 1118.8948
$ws.Range("J20").Value = 1383.375
$ws.Range("K20").Value = 1118.8948
$ws.Range("L20").Value = 1383.375
$ws.Range("M20").Value = -871.8948
$ws.Range("N20").Value = -1877.375

$ws.Range("H94").Value = 768.6889
$ws.Range("I94").Value = 579.17145
$ws.Range("J94").Value = 1432
$ws.Range("K94").Value = 579.17145
$ws.Range("L94").Value = 1432
$ws.Range("M94").Value = -128.17145
$ws.Range("N94").Value = -2334

$ws.Range("H134").Value = 14494779
$ws.Range("I134").Value = 18520242
$ws.Range("J134").Value = 3114.5334
$ws.Range("K134").Value = 55560726
$ws.Range("L134").Value = 9343.600199999999
$ws.Range("M134").Value = -55558191
$ws.Range("N134").Value = -14413.6002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1172.3939
$ws.Range("I31").Value = 856.1
$ws.Range("J31").Value = 4335.3335
$ws.Range("K31").Value = 856.1
$ws.Range("L31").Value = 4335.3335
$ws.Range("M31").Value = -561.1
$ws.Range("N31").Value = -4925.3335

$ws.Range("H34").Value = 1172.3939
$ws.Range("I34").Value = 856.1
$ws.Range("J34").Value = 4335.3335
$ws.Range("K34").Value = 856.1
$ws.Range("L34").Value = 4335.3335
$ws.Range("M34").Value = -654.1
$ws.Range("N34").Value = -4739.3335

$ws.Range("H58").Value = 1190.7142
$ws.Range("I58").Value = 809.65717
$ws.Range("J58").Value = 2143.3572
$ws.Range("K58").Value = 809.65717
$ws.Range("L58").Value = 2143.3572
$ws.Range("M58").Value = -606.65717
$ws.Range("N58").Value = -2549.3572

$ws.Range("H132").Value = 2071.92
$ws.Range("I132").Value = 1892.6666
$ws.Range("J132").Value = 2707.4546
$ws.Range("K132").Value = 5677.9998
$ws.Range("L132").Value = 8122.3638
$ws.Range("M132").Value = -3147.9998
$ws.Range("N132").Value = -13182.3638

$ws.Range("H134").Value = 2519.195
$ws.Range("I134").Value = 1744.8572
$ws.Range("J134").Value = 4187
$ws.Range("K134").Value = 5234.571599999999
$ws.Range("L134").Value = 12561
$ws.Range("M134").Value = -2699.571599999999
$ws.Range("N134").Value = -17631

$ws.Range("H136").Value = 1190.7142
$ws.Range("I136").Value = 809.65717
$ws.Range("J136").Value = 2143.3572
$ws.Range("K136").Value = 2428.97151
$ws.Range("L136").Value = 6430.071599999999
$ws.Range("M136").Value = 121.0284900000001
$ws.Range("N136").Value = -11530.0716

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 9525900
$ws.Range("I131").Value = 604.1429000000001
$ws.Range("J131").Value = 11907224
$ws.Range("K131").Value = 1812.4287
$ws.Range("L131").Value = 35721672
$ws.Range("M131").Value = 3227.5713
$ws.Range("N131").Value = -35731752

$ws.Range("H139").Value = 2870.2856
$ws.Range("I139").Value = 2371.3333
$ws.Range("J139").Value = 3446
$ws.Range("K139").Value = 7113.999899999999
$ws.Range("L139").Value = 10338
$ws.Range("M139").Value = -1973.999899999999
$ws.Range("N139").Value = -20618

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2218.2964
$ws.Range("I102").Value = 2415.8823
$ws.Range("K102").Value = 2415.8823
$ws.Range("M102").Value = -793.8823000000002

$ws.Range("H126").Value = 2184.606
$ws.Range("I126").Value = 1889.2
$ws.Range("J126").Value = 2313.0435
$ws.Range("K126").Value = 5667.6
$ws.Range("L126").Value = 6939.130500000001
$ws.Range("M126").Value = -3197.6
$ws.Range("N126").Value = -11879.1305

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 13461.25
$ws.Range("I22").Value = 1175
$ws.Range("J22").Value = 25747.5
$ws.Range("K22").Value = 1175
$ws.Range("L22").Value = 25747.5
$ws.Range("M22").Value = -880
$ws.Range("N22").Value = -26337.5

$ws.Range("H27").Value = 13461.25
$ws.Range("I27").Value = 1175
$ws.Range("J27").Value = 25747.5
$ws.Range("K27").Value = 1175
$ws.Range("L27").Value = 25747.5
$ws.Range("M27").Value = -1068
$ws.Range("N27").Value = -25961.5

$ws.Range("H61").Value = 1967.2069
$ws.Range("I61").Value = 1987.6522
$ws.Range("J61").Value = 1888.8334
$ws.Range("K61").Value = 1987.6522
$ws.Range("L61").Value = 1888.8334
$ws.Range("M61").Value = -1785.6522
$ws.Range("N61").Value = -2292.8334

$ws.Range("H106").Value = 20507.777
$ws.Range("J106").Value = 20507.777
$ws.Range("L106").Value = 20507.777
$ws.Range("N106").Value = -23031.777

$ws.Range("H113").Value = 1967.2069
$ws.Range("I113").Value = 1987.6522
$ws.Range("J113").Value = 1888.8334
$ws.Range("K113").Value = 1987.6522
$ws.Range("L113").Value = 1888.8334
$ws.Range("M113").Value = 182.3478
$ws.Range("N113").Value = -6228.8334

$ws.Range("H136").Value = 3871.5
$ws.Range("I136").Value = 2357.2258
$ws.Range("J136").Value = 13260
$ws.Range("K136").Value = 7071.6774
$ws.Range("L136").Value = 39780
$ws.Range("M136").Value = -4521.6774
$ws.Range("N136").Value = -44880

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H7").Value = 10819.2
$ws.Range("I7").Value = 10819.2
$ws.Range("K7").Value = 10819.2
$ws.Range("M7").Value = -10706.2

$ws.Range("H17").Value = 2749.25
$ws.Range("I17").Value = 2498.5
$ws.Range("K17").Value = 2498.5
$ws.Range("M17").Value = -2326.5

$ws.Range("H132").Value = 4091.1707
$ws.Range("I132").Value = 4576.207
$ws.Range("J132").Value = 2919
$ws.Range("K132").Value = 13728.621
$ws.Range("L132").Value = 8757
$ws.Range("M132").Value = -11198.621
$ws.Range("N132").Value = -13817
